# Replace the invoice data in row 2 (A2:F2) with new values.
# Several of these look like numbers (leading zeros, long digit strings)
# so we force Text format before assigning, then clear the formatting
# again afterwards so the cells keep their original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:F2")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "00100004010000032593"
$ws.Range("B2").Value = "2024-08-09T07:20:26"
$ws.Range("C2").Value = "3101775072"
$ws.Range("D2").Value = "02"
$ws.Range("E2").Value = "Credito 17"
$ws.Range("F2").Value = "03"

$dataRange.ClearFormats()
